$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - update B2, C2 (D2 unchanged)
$ws.Range("B2").Value = 3519602159820190
$ws.Range("C2").Value = 3519602159820190

# Row 3: RandomForestRegressor - update B3, C3, D3
$ws.Range("B3").Value = 0.02239068914858554
$ws.Range("C3").Value = 0.02378663569108227
$ws.Range("D3").Value = 2354497445318.094

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update B4, C4, D4
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02540866155419177
$ws.Range("C4").Value = 0.02653366009954966
$ws.Range("D4").Value = 0.1150029293793559

# Row 5: AdaBoostRegressor -> MLPRegressor, update B5, C5, D5
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 21143931516222
$ws.Range("C5").Value = 7630438543912.923
$ws.Range("D5").Value = 107916063106830.1
